$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# "Sleep-Funktion" row: points awarded went from 0 to 3 (out of max C22=3)
$ws.Range("D22").Value = 3

# Reviewer comment explaining the partial deduction (new shared string, column F)
$ws.Range("F22").Value = "Hält sich noch nicht ganz an die Beschreibung (WakeUps sind manuell gesetzt, die RESETs bewirken irgendwie so gut wie nichts…)"

# Scroll the window down a bit and move the selection, matching where the
# reviewer was working when they saved the file.
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("E24").Select()
